# "corrección bisección 2.0 (lectura de potencias)"
# Update the bisection-method results table: rows 2-6 get new xi / f(xi) / Error
# values (the bisection now converges toward a different root, reflecting the
# fixed reading of powers), and the former row 7 (iteration "6") is removed
# since convergence is now reached one iteration earlier.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    # Prefixing with an apostrophe forces Excel to store the value as text
    # (matching the workbook's original inline-string cells) instead of
    # inferring a number from content such as "14.9990234375" or
    # "2.40017339336873e-10". Resetting the style back to "Normal" afterwards
    # clears the quote-prefix formatting that the text-entry leaves behind,
    # so the cell keeps its original (unstyled) appearance.
    $ws.Range($Address).Value = "'" + $Text
    $ws.Range($Address).Style = "Normal"
}

# Row 2 (Iteración 1) - xi (B2) is unchanged
Set-TextValue "C2" "14.9990234375"
Set-TextValue "D2" "2.14125859098224"

# Row 3 (Iteración 2)
Set-TextValue "B3" "-2.85874140901776"
Set-TextValue "C3" "-1.8466024818928"
Set-TextValue "D3" "0.289659790013761"

# Row 4 (Iteración 3)
Set-TextValue "B4" "-3.14840119903152"
Set-TextValue "C4" "-0.100289494253238"
Set-TextValue "D4" "0.0158040433021078"

# Row 5 (Iteración 4)
Set-TextValue "B5" "-3.16420524233363"
Set-TextValue "C5" "-0.0002491447552444"
Set-TextValue "D5" "3.92617396967054e-05"

# Row 6 (Iteración 5)
Set-TextValue "B6" "-3.16424450407332"
Set-TextValue "C6" "-1.52308840507432e-09"
Set-TextValue "D6" "2.40017339336873e-10"

# The old row 7 (Iteración 6) is no longer needed - bisection now converges
# one step sooner, so drop the last row entirely and let the sheet
# dimension shrink from A1:D7 to A1:D6 automatically.
$ws.Rows("7:7").Delete()
